# Add the new "2022-Q4" sheet right after "总计" and before "2022-Q3", and
# insert the corresponding summary row into the "总计" sheet.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item(1)
$q3Sheet = $wb.Worksheets.Item(2)

# --- 1. Create the new "2022-Q4" worksheet right after "总计" -------------
# Duplicate the existing "2022-Q3" sheet (same column layout/styling) and
# drop it in right after "总计"; then overwrite its contents.
$q3Sheet.Copy($null, $totalSheet)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

# The source sheet only has 7 rows (header + 6 data rows); the new sheet
# needs 9 (header + 8 data rows), so insert 2 more data rows first.
$newSheet.Rows.Item(8).Insert()
$newSheet.Rows.Item(8).Insert()

$q4data = @(
    @(0, "012367", "上投摩根安荣回报混合C", "12.39", "25.70", "0.95", "0.1177", 10),
    @(1, "004738", "上投摩根安隆回报混合A", "10.31", "23.32", "0.97", "0.1000", 6),
    @(2, "012366", "上投摩根安荣回报混合A", "9.72",  "25.70", "0.95", "0.0923", 10),
    @(3, "004823", "上投摩根安裕回报混合A", "4.26",  "36.12", "1.66", "0.0707", 8),
    @(4, "004824", "上投摩根安裕回报混合C", "3.64",  "36.12", "1.66", "0.0604", 8),
    @(5, "004739", "上投摩根安隆回报混合C", "4.60",  "23.32", "0.97", "0.0446", 6),
    @(6, "001231", "银华泰利灵活配置混合A", "0.87",  "22.82", "1.08", "0.0094", 3),
    @(7, "002328", "银华泰利灵活配置混合C", "0.03",  "22.82", "1.08", "0.0003", 3)
)

$r = 2
foreach ($row in $q4data) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]

    $textRange = $newSheet.Range($newSheet.Cells.Item($r, 2), $newSheet.Cells.Item($r, 7))
    $textRange.NumberFormat = "@"
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Value = $row[6]
    $textRange.ClearFormats()

    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Re-apply the column-A style to the two freshly inserted rows (8 & 9) so
# they match the rest of the A column (bold/centered with border).
$newSheet.Range("A2").Copy()
$newSheet.Range("A8:A9").PasteSpecial(-4122)
$newSheet.Range("A1").Select()

# --- 2. Insert the 2022-Q4 summary row into "总计" --------------------------
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("A2:D2").ClearFormats()
$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q4"
$totalSheet.Cells.Item(2, 3).Value = 8
$totalSheet.Cells.Item(2, 4).Value = 0.5

# Re-apply the column-A style (lost on row insert) and fix the sequential
# index numbers (A2..A6 should read 0,1,2,3,4) now that every row moved
# down by one.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(5, 1).Value = 3
$totalSheet.Cells.Item(6, 1).Value = 4

$totalSheet.Range("A1").Select()
